{"js": "// Add sleuth-zipkin for logs:\n// 1. Highlight the \"Config\" list item in yellow (paragraph mark + run).\n// 2. Remove the \"Gateway\" list item entirely.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet gatewayParagraph = null;\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text.trim();\n  if (text === \"Config\") {\n    paragraph.font.highlightColor = \"Yellow\";\n  } else if (text === \"Gateway\") {\n    gatewayParagraph = paragraph;\n  }\n}\n\nif (gatewayParagraph) {\n  gatewayParagraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Add sleuth-zipkin for logs:\n# 1. Highlight the \"Config\" list item in yellow (paragraph mark + run).\n# 2. Remove the \"Gateway\" list item entirely.\n\n$d = $word.ActiveDocument\n\n# Snapshot paragraphs first since we mutate (delete) the collection below.\n$paragraphs = @($d.Paragraphs)\n\nforeach ($p in $paragraphs) {\n    $text = $p.Range.Text.Trim()\n    if ($text -eq \"Config\") {\n        # wdYellow = 7; applying to Range.Font (not just Range) also stamps\n        # the paragraph mark's run properties, matching Word's own behavior\n        # when highlighting a whole list-item line.\n        $p.Range.Font.HighlightColorIndex = 7\n    } elseif ($text -eq \"Gateway\") {\n        $p.Range.Delete()\n    }\n}\n"}
